# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-level detail) right before the
#    "总计" (summary) sheet.
# 2. Insert a new summary row for "2022-Q1" at the top of the "总计" sheet's
#    data (shifting the other quarters down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a value to be written as literal TEXT (not auto-coerced to a
# number by Excel) while leaving the cell's style untouched (style "Normal").
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q1" worksheet
#
# NOTE: sheet object references in this host resolve by position, so a
# reference captured for a sheet AT OR AFTER the insertion point (here,
# "总计") goes stale once a new sheet is spliced in front of it. We grab
# "2021-Q4" (unaffected, since it sits before the insertion point) up front,
# but we re-fetch "总计" by name *after* the insert+rename are both done.
# ---------------------------------------------------------------------------
$q4sheet = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Header row (B1:H1) - copy style from an existing header cell ("s=2") then
# overwrite the text.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $srcHeader = $q4sheet.Cells.Item(1, $col)
    $dstHeader = $newSheet.Cells.Item(1, $col)
    $srcHeader.Copy($dstHeader)
    $dstHeader.Value = $headers[$i]
}

# Data rows
$rows = @(
    @("001302", "前海开源金银珠宝主题精选混合A",         "8.61",  "91.91", "7.78", "0.6699", 9),
    @("003304", "前海开源沪港深核心资源灵活配置混合A",     "5.91",  "93.10", "7.60", "0.4492", 6),
    @("001468", "广发改革先锋灵活配置混合",               "10.51", "76.65", "3.14", "0.3300", 4),
    @("004475", "华泰柏瑞富利灵活配置混合",               "4.43",  "93.22", "6.62", "0.2933", 3),
    @("002207", "前海开源金银珠宝主题精选混合C",         "3.45",  "91.91", "7.78", "0.2684", 9),
    @("003305", "前海开源沪港深核心资源灵活配置混合C",     "2.19",  "93.10", "7.60", "0.1664", 6),
    @("001247", "华泰柏瑞新利灵活配置混合A",             "10.68", "20.70", "1.15", "0.1228", 3),
    @("004010", "华泰柏瑞鼎利灵活配置混合A",             "10.24", "21.26", "1.19", "0.1219", 3),
    @("002091", "华泰柏瑞新利灵活配置混合C",             "5.98",  "20.70", "1.15", "0.0688", 3),
    @("004011", "华泰柏瑞鼎利灵活配置混合C",             "3.47",  "21.26", "1.19", "0.0413", 3)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = 2 + $r
    $data = $rows[$r]

    # Column A: plain 0-based index, styled like the header ("s=2")
    $srcA = $q4sheet.Cells.Item(2, 1)
    $dstA = $newSheet.Cells.Item($rowNum, 1)
    $srcA.Copy($dstA)
    $dstA.Value = $r

    # Columns B-G: text values (fund code / name / scale / position figures)
    Set-TextValue $newSheet.Cells.Item($rowNum, 2) $data[0]
    Set-TextValue $newSheet.Cells.Item($rowNum, 3) $data[1]
    Set-TextValue $newSheet.Cells.Item($rowNum, 4) $data[2]
    Set-TextValue $newSheet.Cells.Item($rowNum, 5) $data[3]
    Set-TextValue $newSheet.Cells.Item($rowNum, 6) $data[4]
    Set-TextValue $newSheet.Cells.Item($rowNum, 7) $data[5]

    # Column H: numeric rank
    $newSheet.Cells.Item($rowNum, 8).Value = $data[6]
}

# ---------------------------------------------------------------------------
# 2. Add the "2022-Q1" row to the "总计" summary sheet (new top data row)
# ---------------------------------------------------------------------------
# Re-fetch fresh now that the sheet list has been altered above.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.UsedRange.ClearContents()

$totalSheet.Cells.Item(1, 2).Value = "日期"
$totalSheet.Cells.Item(1, 3).Value = "持有数量(只)"
$totalSheet.Cells.Item(1, 4).Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 10, 2.53),
    @("2021-Q4", 6, 1.35),
    @("2021-Q3", 5, 1.19),
    @("2021-Q2", 8, 2.69),
    @("2021-Q1", 6, 1.15),
    @("2020-Q4", 4, 2.5)
)

# Row 7 is brand new (sheet used to stop at row 6) - its "A" cell needs the
# same style ("s=2") the other index cells in column A already carry.
$styledA = $totalSheet.Cells.Item(6, 1)
$newA7 = $totalSheet.Cells.Item(7, 1)
$styledA.Copy($newA7)

for ($r = 0; $r -lt $summaryRows.Length; $r++) {
    $rowNum = 2 + $r
    $data = $summaryRows[$r]

    $totalSheet.Cells.Item($rowNum, 1).Value = $r
    $totalSheet.Cells.Item($rowNum, 2).Value = $data[0]
    $totalSheet.Cells.Item($rowNum, 3).Value = $data[1]
    $totalSheet.Cells.Item($rowNum, 4).Value = $data[2]
}
